$d = $word.ActiveDocument

$d.Content.Find.Execute("461×9=4149", $true, $false, $false, $false, $false, $true, 1, $false, "618×6=3708", 2)
$d.Content.Find.Execute("244×2=488", $true, $false, $false, $false, $false, $true, 1, $false, "643×9=5787", 2)
$d.Content.Find.Execute("517×3=1551", $true, $false, $false, $false, $false, $true, 1, $false, "670×8=5360", 2)
$d.Content.Find.Execute("153×5=765", $true, $false, $false, $false, $false, $true, 1, $false, "440×5=2200", 2)
$d.Content.Find.Execute("547×2=1094", $true, $false, $false, $false, $false, $true, 1, $false, "778×3=2334", 2)
$d.Content.Find.Execute("682×4=2728", $true, $false, $false, $false, $false, $true, 1, $false, "486×7=3402", 2)
$d.Content.Find.Execute("526×5=2630", $true, $false, $false, $false, $false, $true, 1, $false, "390×3=1170", 2)
$d.Content.Find.Execute("474×5=2370", $true, $false, $false, $false, $false, $true, 1, $false, "541×8=4328", 2)
$d.Content.Find.Execute("238×2=476", $true, $false, $false, $false, $false, $true, 1, $false, "266×6=1596", 2)
$d.Content.Find.Execute("257×3=771", $true, $false, $false, $false, $false, $true, 1, $false, "779×5=3895", 2)
$d.Content.Find.Execute("433×4=1732", $true, $false, $false, $false, $false, $true, 1, $false, "990×3=2970", 2)
$d.Content.Find.Execute("103×9=927", $true, $false, $false, $false, $false, $true, 1, $false, "170×4=680", 2)
$d.Content.Find.Execute("362×9=3258", $true, $false, $false, $false, $false, $true, 1, $false, "740×8=5920", 2)
$d.Content.Find.Execute("808×3=2424", $true, $false, $false, $false, $false, $true, 1, $false, "634×9=5706", 2)
$d.Content.Find.Execute("592×3=1776", $true, $false, $false, $false, $false, $true, 1, $false, "162×7=1134", 2)
$d.Content.Find.Execute("725×7=5075", $true, $false, $false, $false, $false, $true, 1, $false, "243×5=1215", 2)
$d.Content.Find.Execute("483×5=2415", $true, $false, $false, $false, $false, $true, 1, $false, "249×8=1992", 2)
$d.Content.Find.Execute("128×5=640", $true, $false, $false, $false, $false, $true, 1, $false, "471×3=1413", 2)
$d.Content.Find.Execute("366×5=1830", $true, $false, $false, $false, $false, $true, 1, $false, "755×5=3775", 2)
$d.Content.Find.Execute("310×4=1240", $true, $false, $false, $false, $false, $true, 1, $false, "966×6=5796", 2)
$d.Content.Find.Execute("929×3=2787", $true, $false, $false, $false, $false, $true, 1, $false, "722×9=6498", 2)
$d.Content.Find.Execute("423×3=1269", $true, $false, $false, $false, $false, $true, 1, $false, "695×8=5560", 2)
$d.Content.Find.Execute("414×9=3726", $true, $false, $false, $false, $false, $true, 1, $false, "481×2=962", 2)
$d.Content.Find.Execute("328×8=2624", $true, $false, $false, $false, $false, $true, 1, $false, "829×5=4145", 2)
$d.Content.Find.Execute("926×2=1852", $true, $false, $false, $false, $false, $true, 1, $false, "692×3=2076", 2)
